$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - human readable headers (replacing slug-style names)
$ws.Range("A1").Value = "Combustible"
$ws.Range("B1").Value = "Comarca nombre"
$ws.Range("C1").Value = "Número hogares"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "Provincia código"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Provincia nombre"
$ws.Range("H1").Value = "Año de construcción"

# Row 2 - sdmx/iaest measure & dimension identifiers
$ws.Range("A2").Value = "iaest-measure:combustible"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "iaest-measure:numero-hogares"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-dimension:ano-de-construccion"

# Row 3 - medida/dim markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "dim"
$ws.Range("H3").Value = "dim"

# Row 4 - datatype / concept scheme references
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "URI-comarca"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "URI-Provincia"
$ws.Range("H4").Value = "skos:Concept"

# Row 5 - mapping file moved from column B to column H; remove the old
# cell entirely (delete, not just clear, so no stray empty cell remains)
# before writing the new one, then copy the row's standard cell format
# (style index shared by every populated cell in the sheet) onto H5.
$ws.Range("B5").Delete()
$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = "mapping-ano-de-construccion.xlsx"
